$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2:A4").UnMerge()
$ws.Range("A5:A6").ClearContents()
$ws.Range("A2:A6").Merge()
$ws.Range("A2:A6").Select()
